$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Stefano%Amatori%NULL%1,          Sabrina%Donati Zeppa%NULL%2,          Sabrina%Donati Zeppa%NULL%0,          Antonio%Preti%NULL%2,          Antonio%Preti%NULL%0,          Marco%Gervasi%NULL%2,          Marco%Gervasi%NULL%0,          Erica%Gobbi%NULL%2,          Erica%Gobbi%NULL%0,          Fabio%Ferrini%NULL%2,          Fabio%Ferrini%NULL%0,          Marco B. L.%Rocchi%NULL%1,          Carlo%Baldari%NULL%2,          Carlo%Baldari%NULL%0,          Fabrizio%Perroni%NULL%2,          Fabrizio%Perroni%NULL%0,          Giovanni%Piccoli%NULL%2,          Giovanni%Piccoli%NULL%0,          Vilberto%Stocchi%NULL%1,          Piero%Sestili%NULL%1,          Davide%Sisti%NULL%2,          Davide%Sisti%NULL%0]'
$ws.Range('E3').Value = '[Imran%Aslan%NULL%1,          Dominika%Ochnik%NULL%1,          Orhan%Çınar%NULL%1]'
$ws.Range('E5').Value = '[Biber%D%coreGivesNoEmail%1,        March%L%coreGivesNoEmail%1,        Marshall%S.%coreGivesNoEmail%1,        Schmidt%SW%coreGivesNoEmail%1]'
$ws.Range('E6').Value = '[Stéphanie%Bourion-Bédès%NULL%1,          Cyril%Tarquinio%NULL%1,          Martine%Batt%NULL%1,          Pascale%Tarquinio%NULL%1,          Romain%Lebreuilly%NULL%1,          Christine%Sorsana%NULL%1,          Karine%Legrand%NULL%1,          Hélène%Rousseau%NULL%1,          Cédric%Baumann%NULL%1]'
$ws.Range('E7').Value = '[Rong-ning%Chen%NULL%1,          Shun-wei%Liang%NULL%1,          Yang%Peng%NULL%1,          Xue-guo%Li%NULL%1,          Jian-bin%Chen%NULL%1,          Si-yao%Tang%NULL%1,          Jing-bo%Zhao%NULL%1]'
$ws.Range('E8').Value = '[Xinli%Chi%NULL%0,          Benjamin%Becker%NULL%0,          Qian%Yu%NULL%0,          Peter%Willeit%NULL%0,          Can%Jiao%NULL%0,          Liuyue%Huang%NULL%0,          M. Mahhub%Hossain%NULL%0,          Igor%Grabovac%NULL%0,          Albert%Yeung%NULL%0,          Jingyuan%Lin%NULL%0,          Nicola%Veronese%NULL%0,          Jian%Wang%NULL%0,          Xinqi%Zhou%NULL%0,          Scott R.%Doig%NULL%0,          Xiaofeng%Liu%NULL%0,          Andre F.%Carvalho%NULL%0,          Lin%Yang%NULL%0,          Tao%Xiao%NULL%0,          Liye%Zou%NULL%0,          Paolo%Fusar-Poli%NULL%0,          Marco%Solmi%NULL%0]'
$ws.Range('E9').Value = '[Bablu Kumar%Dhar%drbablu.ba@diu.edu.bd%1,          Foster Kofi%Ayittey%NULL%2,          Foster Kofi%Ayittey%NULL%0,          Sabrina Maria%Sarkar%NULL%1]'
$ws.Range('E10').Value = '[Liu%X.%coreGivesNoEmail%1,        Lovibond%P.%coreGivesNoEmail%1,        Lozano%M.%coreGivesNoEmail%1,        Mosqueda%D\u00edaz A.%coreGivesNoEmail%1,        Ozamiz-Etxebarria%N.%coreGivesNoEmail%1]'
$ws.Range('E11').Value = '[Julia%Dratva%NULL%1,          Annina%Zysset%NULL%2,          Annina%Zysset%NULL%0,          Nadine%Schlatter%NULL%1,          Agnes%von Wyl%NULL%1,          Marion%Huber%NULL%1,          Thomas%Volken%NULL%2,          Thomas%Volken%NULL%0]'
$ws.Range('E12').Value = '[Chen%Du%NULL%1,          Megan Chong Hueh%Zan%NULL%2,          Megan Chong Hueh%Zan%NULL%0,          Min Jung%Cho%NULL%1,          Jenifer I.%Fenton%NULL%2,          Jenifer I.%Fenton%NULL%0,          Pao Ying%Hsiao%NULL%2,          Pao Ying%Hsiao%NULL%0,          Richard%Hsiao%NULL%1,          Laura%Keaver%NULL%1,          Chang-Chi%Lai%NULL%2,          Chang-Chi%Lai%NULL%0,          HeeSoon%Lee%NULL%1,          Mary-Jon%Ludy%NULL%2,          Mary-Jon%Ludy%NULL%0,          Wan%Shen%NULL%1,          Winnie Chee Siew%Swee%NULL%2,          Winnie Chee Siew%Swee%NULL%0,          Jyothi%Thrivikraman%NULL%2,          Jyothi%Thrivikraman%NULL%0,          Kuo-Wei%Tseng%NULL%1,          Wei-Chin%Tseng%NULL%1,          Robin M.%Tucker%NULL%1]'
$ws.Range('E13').Value = '[Aziz%Essadek%NULL%1,          Thomas%Rabeyron%NULL%1]'
$ws.Range('C14').Value = '"E\u2010learning: Depression, anxiety, and stress symptomatology among Lebanese university students during COVID\u201019 quarantine"'
$ws.Range('E14').Value = '[Mirna%Fawaz%xref no email%1,  Ali%Samaha%xref no email%1]'
$ws.Range('F14').Value = '10.1111/nuf.12521'
$ws.Range('G14').Value = 'CROSSREF'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '2023-05-25'
$ws.Range('H14').ClearFormats()
$ws.Range('E15').Value = '[Wenning%Fu%NULL%1,          Shijiao%Yan%NULL%0,          Qiao%Zong%NULL%1,          Dan%Anderson-Luxford%NULL%1,          Xingyue%Song%NULL%1,          Zhiyue%Lv%NULL%1,          Chuanzhu%Lv%NULL%0]'
$ws.Range('E16').Value = '[Beata%Gavurova%NULL%1,          Viera%Ivankova%NULL%2,          Viera%Ivankova%NULL%0,          Martin%Rigelsky%NULL%2,          Martin%Rigelsky%NULL%0]'
$ws.Range('E17').Value = '[Eman R%Ghazawy%emanghazawy@yahoo.com%1,          Ashraf A%Ewis%NULL%1,          Eman M%Mahfouz%NULL%1,          Doaa M%Khalil%NULL%1,          Ahmed%Arafa%NULL%1,          Zeinab%Mohammed%NULL%1,          El-Nabgha F%Mohammed%NULL%1,          Ebtesam E%Hassan%NULL%1,          Sarah%Abdel Hamid%NULL%1,          Somaya A%Ewis%NULL%1,          Abd El-Nassir S%Mohammed%NULL%1]'
